# Update progress values for a handful of clinical trials in TrialsSetup.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B5").Value  = 88   # KATALYST
$ws.Range("B7").Value  = 38   # CLEAR
$ws.Range("B8").Value  = 12   # CADANCE
$ws.Range("B10").Value = 12   # OPERA-2

$wb.Save()
